# Apply the recorded edits to rows 16-18 of the "Artfynd" sheet.
#
# Summary of the change:
#  - Row 16: Q16/R16 (Ost/Nord coordinates) are rounded from long decimals
#    down to plain integers; the Starttid/Sluttid cells (Z16/AB16, both
#    "00:00") are cleared out entirely.
#  - Rows 17 and 18 swap their entire record content (species id, taxon,
#    names, author, coordinates, and the Age/Sex/Activity/Method columns
#    K/L/M/N that only the "Spillkråka" / Dryocopus martius record carries),
#    with the Ost/Nord coordinates rounded to integers in their new spot.
#  - Rows 17 and 18 also lose their Starttid/Sluttid cells, same as row 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: round Ost/Nord, drop Starttid/Sluttid -------------------------
$ws.Range("Q16").Value = 653207
$ws.Range("R16").Value = 6600027
$ws.Range("Z16").ClearContents()
$ws.Range("AB16").ClearContents()

# --- Capture the current (pre-edit) row 17 and row 18 record data ---------
$r17_A = $ws.Range("A17").Value2
$r17_B = $ws.Range("B17").Value2
$r17_D = $ws.Range("D17").Value2
$r17_E = $ws.Range("E17").Value2
$r17_F = $ws.Range("F17").Value2
$r17_G = $ws.Range("G17").Value2
$r17_H = $ws.Range("H17").Value2

$r18_A = $ws.Range("A18").Value2
$r18_B = $ws.Range("B18").Value2
$r18_D = $ws.Range("D18").Value2
$r18_E = $ws.Range("E18").Value2
$r18_F = $ws.Range("F18").Value2
$r18_G = $ws.Range("G18").Value2
$r18_H = $ws.Range("H18").Value2
$r18_M = $ws.Range("M18").Value2

# --- Row 17 becomes the old row 18 record (Spillkråka / Dryocopus martius) -
$ws.Range("A17").Value = $r18_A
$ws.Range("B17").Value = $r18_B
$ws.Range("D17").Value = $r18_D
$ws.Range("E17").Value = $r18_E
$ws.Range("F17").Value = $r18_F
$ws.Range("G17").Value = $r18_G
$ws.Range("H17").Value = $r18_H
# K18/L18/N18 were present-but-empty inline-string cells; materialize the
# same empty cells at K17/L17/N17 (plain Value="" does not create a cell).
$ws.Range("I17").Copy($ws.Range("K17"))
$ws.Range("I17").Copy($ws.Range("L17"))
$ws.Range("M17").Value = $r18_M
$ws.Range("I17").Copy($ws.Range("N17"))
$ws.Range("Q17").Value = 653148
$ws.Range("R17").Value = 6600341
$ws.Range("Z17").ClearContents()
$ws.Range("AB17").ClearContents()

# --- Row 18 becomes the old row 17 record (Blåsippa / Hepatica nobilis) ---
$ws.Range("A18").Value = $r17_A
$ws.Range("B18").Value = $r17_B
$ws.Range("D18").Value = $r17_D
$ws.Range("E18").Value = $r17_E
$ws.Range("F18").Value = $r17_F
$ws.Range("G18").Value = $r17_G
$ws.Range("H18").Value = $r17_H
$ws.Range("K18").ClearContents()
$ws.Range("L18").ClearContents()
$ws.Range("M18").ClearContents()
$ws.Range("N18").ClearContents()
$ws.Range("Q18").Value = 653206
$ws.Range("R18").Value = 6599944
$ws.Range("Z18").ClearContents()
$ws.Range("AB18").ClearContents()
